$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")
$ws.Rows.Item(25).Insert()
$ws.Range("A25").Copy()
$ws.Range("D25").PasteSpecial(-4122)
